$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1038.1428
$ws.Range("I6").Value = 253.5
$ws.Range("K6").Value = 760.5
$ws.Range("M6").Value = -648.5
$ws.Range("H11").Value = 26.25
$ws.Range("I11").Value = 26.25
$ws.Range("K11").Value = 26.25
$ws.Range("M11").Value = 113.75
$ws.Range("H15").Value = 1845.1464
$ws.Range("I15").Value = 1845.1464
$ws.Range("K15").Value = 5535.439200000001
$ws.Range("M15").Value = -5366.439200000001
$ws.Range("H17").Value = 2701.3438
$ws.Range("J17").Value = 2773.8965
$ws.Range("L17").Value = 8321.6895
$ws.Range("N17").Value = -8657.6895
$ws.Range("H32").Value = 1783.3334
$ws.Range("J32").Value = 2275
$ws.Range("L32").Value = 2275
$ws.Range("N32").Value = -2927
$ws.Range("H40").Value = 4156.269
$ws.Range("I40").Value = 3557
$ws.Range("K40").Value = 3557
$ws.Range("M40").Value = -3382
$ws.Range("H74").Value = 7499.8
$ws.Range("I74").Value = 6499.6665
$ws.Range("K74").Value = 6499.6665
$ws.Range("M74").Value = -5563.6665
$ws.Range("H76").Value = 4249.5
$ws.Range("I76").Value = 6999
$ws.Range("J76").Value = 1500
$ws.Range("K76").Value = 6999
$ws.Range("L76").Value = 1500
$ws.Range("M76").Value = -6684
$ws.Range("N76").Value = -2130
$ws.Range("H77").Value = 7499.8
$ws.Range("I77").Value = 6499.6665
$ws.Range("K77").Value = 32498.3325
$ws.Range("M77").Value = -27818.3325
$ws.Range("H79").Value = 4249.5
$ws.Range("I79").Value = 6999
$ws.Range("J79").Value = 1500
$ws.Range("K79").Value = 6999
$ws.Range("L79").Value = 1500
$ws.Range("M79").Value = -5907
$ws.Range("N79").Value = -3684
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -39992
$ws.Range("H107").Value = 141
$ws.Range("I107").Value = 149.5
$ws.Range("J107").Value = 90
$ws.Range("K107").Value = 149.5
$ws.Range("L107").Value = 90
$ws.Range("M107").Value = 1770.5
$ws.Range("N107").Value = -3930
$ws.Range("H111").Value = 945.6
$ws.Range("I111").Value = 686.7778
$ws.Range("J111").Value = 3275
$ws.Range("K111").Value = 2060.3334
$ws.Range("L111").Value = 9825
$ws.Range("M111").Value = 1006.6666
$ws.Range("N111").Value = -15959
$ws.Range("H132").Value = 5234.154
$ws.Range("I132").Value = 1640.3636
$ws.Range("K132").Value = 4921.0908
$ws.Range("M132").Value = -2391.0908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 43.333332
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 50
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 62
$ws.Range("N5").Value = -254
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = ""
$ws.Range("H92").Value = 35000
$ws.Range("J92").Value = 35000
$ws.Range("L92").Value = 35000
$ws.Range("N92").Value = -39992
$ws.Range("H94").Value = 45000
$ws.Range("J94").Value = 45000
$ws.Range("L94").Value = 45000
$ws.Range("N94").Value = -46802
$ws.Range("H110").Value = 3903.1428
$ws.Range("I110").Value = 3164.4
$ws.Range("K110").Value = 3164.4
$ws.Range("M110").Value = -1119.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 43.333332
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = 65
$ws.Range("N4").Value = -260
$ws.Range("H22").Value = 175
$ws.Range("I22").Value = 150
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 150
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = 23
$ws.Range("N22").Value = -596
$ws.Range("H80").Value = 169
$ws.Range("I80").Value = 124.8
$ws.Range("J80").Value = 187.41667
$ws.Range("K80").Value = 124.8
$ws.Range("L80").Value = 187.41667
$ws.Range("M80").Value = 873.2
$ws.Range("N80").Value = -2183.41667
$ws.Range("H83").Value = 169
$ws.Range("I83").Value = 124.8
$ws.Range("J83").Value = 187.41667
$ws.Range("K83").Value = 624
$ws.Range("L83").Value = 937.0833500000001
$ws.Range("M83").Value = 4368
$ws.Range("N83").Value = -10921.08335
$ws.Range("H86").Value = 4640.6924
$ws.Range("J86").Value = 6000.5713
$ws.Range("L86").Value = 6000.5713
$ws.Range("N86").Value = -8246.5713
$ws.Range("H89").Value = 4640.6924
$ws.Range("J89").Value = 6000.5713
$ws.Range("L89").Value = 30002.8565
$ws.Range("N89").Value = -41234.85649999999
$ws.Range("H94").Value = 325.5
$ws.Range("I94").Value = 264.18182
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 264.18182
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = 186.81818
$ws.Range("N94").Value = -1902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 59.666668
$ws.Range("I7").Value = 59.666668
$ws.Range("K7").Value = 59.666668
$ws.Range("M7").Value = 53.333332
$ws.Range("H58").Value = 4344.5
$ws.Range("I58").Value = 1362.4
$ws.Range("K58").Value = 1362.4
$ws.Range("M58").Value = -1159.4
$ws.Range("H99").Value = 3999.8
$ws.Range("H102").Value = 36675
$ws.Range("J102").Value = 36675
$ws.Range("L102").Value = 36675
$ws.Range("N102").Value = -41543
$ws.Range("H104").Value = 39375
$ws.Range("I104").Value = 39375
$ws.Range("K104").Value = 39375
$ws.Range("M104").Value = -36754
$ws.Range("H107").Value = 526.5
$ws.Range("I107").Value = 552
$ws.Range("J107").Value = 399
$ws.Range("K107").Value = 552
$ws.Range("L107").Value = 399
$ws.Range("M107").Value = 1368
$ws.Range("N107").Value = -4239
$ws.Range("H126").Value = 3999.8
$ws.Range("H132").Value = 4749.154
$ws.Range("I132").Value = 4110.5557
$ws.Range("K132").Value = 12331.6671
$ws.Range("M132").Value = -9801.667099999999
$ws.Range("H136").Value = 4344.5
$ws.Range("I136").Value = 1362.4
$ws.Range("K136").Value = 4087.2
$ws.Range("M136").Value = -1537.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 15335060
$ws.Range("I7").Value = 13750775
$ws.Range("K7").Value = 13750775
$ws.Range("M7").Value = -13750663
$ws.Range("H8").Value = 15335060
$ws.Range("I8").Value = 13750775
$ws.Range("K8").Value = 13750775
$ws.Range("M8").Value = -13750636
$ws.Range("H11").Value = 10690229
$ws.Range("I11").Value = 8095796.5
$ws.Range("K11").Value = 8095796.5
$ws.Range("M11").Value = -8095657.5
$ws.Range("H15").Value = 37554.445
$ws.Range("J15").Value = 37554.445
$ws.Range("L15").Value = 37554.445
$ws.Range("N15").Value = -38130.445
$ws.Range("H81").Value = 37554.445
$ws.Range("J81").Value = 37554.445
$ws.Range("L81").Value = 37554.445
$ws.Range("N81").Value = -39550.445
$ws.Range("H84").Value = 37554.445
$ws.Range("J84").Value = 37554.445
$ws.Range("L84").Value = 112663.335
$ws.Range("N84").Value = -122647.335
$ws.Range("H107").Value = 265
$ws.Range("I107").Value = 215.71428
$ws.Range("K107").Value = 215.71428
$ws.Range("M107").Value = 1704.28572
$ws.Range("H113").Value = 7365.077
$ws.Range("I113").Value = 1687.25
$ws.Range("K113").Value = 1687.25
$ws.Range("M113").Value = 482.75
$ws.Range("H122").Value = 2678.3333
$ws.Range("J122").Value = 4198.6
$ws.Range("L122").Value = 12595.8
$ws.Range("N122").Value = -17495.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5196.524
$ws.Range("I46").Value = 3512.7
$ws.Range("K46").Value = 3512.7
$ws.Range("M46").Value = -3324.7
$ws.Range("H61").Value = 3543.6316
$ws.Range("I61").Value = 1444.25
$ws.Range("K61").Value = 1444.25
$ws.Range("M61").Value = -1242.25
$ws.Range("H113").Value = 3543.6316
$ws.Range("I113").Value = 1444.25
$ws.Range("K113").Value = 1444.25
$ws.Range("M113").Value = 725.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 837.625
$ws.Range("I113").Value = 928.7143
$ws.Range("K113").Value = 2786.1429
$ws.Range("M113").Value = -616.1428999999998
